$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.200908310620848
$ws.Range("C2").Value = 0.200908310620848
$ws.Range("D2").Value = 0.190099000331555
$ws.Range("E2").Value = 0.00125887971024484
$ws.Range("F2").Value = 0.6722

# Row 3
$ws.Range("B3").Value = 6.83085533184107
$ws.Range("C3").Value = 6.83085533184107
$ws.Range("D3").Value = 6.46334024699978
$ws.Range("E3").Value = 0.0428017395313272
$ws.Range("F3").Value = 0.0093

# Row 4
$ws.Range("B4").Value = 0.373120665870231
$ws.Range("C4").Value = 0.373120665870231
$ws.Range("D4").Value = 0.353045950990217
$ws.Range("E4").Value = 0.00233795224441222
$ws.Range("F4").Value = 0.5579

# Row 5
$ws.Range("B5").Value = 152.18805295632
$ws.Range("C5").Value = 1.05686147886333
$ws.Range("E5").Value = 0.953601428514016

# Row 6
$ws.Range("B6").Value = 159.592937264652
